# Updated cryptos list on Thu Mar  2 13:40:29 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and renames row 51 from Flow to EOS (coin/link/price/volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text, and whether the text must be
# "number-protected" (it looks like a number, e.g. "1.001" or "21.97",
# so a plain .Value assignment would be auto-coerced to a numeric value
# and lose the exact formatting -- e.g. "1.000" would become 1). For those
# cells we temporarily force a Text number format, assign the literal
# string, then clear the format again so the cell keeps its original
# (unstyled) appearance while the stored value remains the exact text.
$updates = @(
    @{ Cell = "D2"; Value = '23.389.14'; Protect = $false }
    @{ Cell = "E2"; Value = '  -1.30%  '; Protect = $false }
    @{ Cell = "D3"; Value = '1.636.38'; Protect = $false }
    @{ Cell = "E3"; Value = '  -0.98%  '; Protect = $false }
    @{ Cell = "D4"; Value = '1.001'; Protect = $true }
    @{ Cell = "E4"; Value = '  -0.16%  '; Protect = $false }
    @{ Cell = "D5"; Value = '1.001'; Protect = $true }
    @{ Cell = "E5"; Value = '  -0.09%  '; Protect = $false }
    @{ Cell = "D6"; Value = '299.69'; Protect = $true }
    @{ Cell = "E6"; Value = '  -1.17%  '; Protect = $false }
    @{ Cell = "D7"; Value = '0.3783'; Protect = $true }
    @{ Cell = "E7"; Value = '  -0.49%  '; Protect = $false }
    @{ Cell = "D8"; Value = '50.33'; Protect = $true }
    @{ Cell = "E8"; Value = '  -1.48%  '; Protect = $false }
    @{ Cell = "D9"; Value = '0.3525'; Protect = $true }
    @{ Cell = "E9"; Value = '  -2.57%  '; Protect = $false }
    @{ Cell = "D10"; Value = '0.08054'; Protect = $true }
    @{ Cell = "E10"; Value = '  -1.86%  '; Protect = $false }
    @{ Cell = "D11"; Value = '1.209'; Protect = $true }
    @{ Cell = "E11"; Value = '  -2.90%  '; Protect = $false }
    @{ Cell = "D12"; Value = '1.002'; Protect = $true }
    @{ Cell = "E12"; Value = '  -0.16%  '; Protect = $false }
    @{ Cell = "D13"; Value = '21.97'; Protect = $true }
    @{ Cell = "E13"; Value = '  -2.87%  '; Protect = $false }
    @{ Cell = "D14"; Value = '6.336'; Protect = $true }
    @{ Cell = "E14"; Value = '  -2.77%  '; Protect = $false }
    @{ Cell = "D15"; Value = '7.260'; Protect = $true }
    @{ Cell = "E15"; Value = '  -2.52%  '; Protect = $false }
    @{ Cell = "D16"; Value = '0.00001203'; Protect = $true }
    @{ Cell = "E16"; Value = '  -2.40%  '; Protect = $false }
    @{ Cell = "D17"; Value = '1.636.17'; Protect = $false }
    @{ Cell = "E17"; Value = '  -0.83%  '; Protect = $false }
    @{ Cell = "D18"; Value = '95.97'; Protect = $true }
    @{ Cell = "E18"; Value = '  -1.38%  '; Protect = $false }
    @{ Cell = "D19"; Value = '0.06952'; Protect = $true }
    @{ Cell = "E19"; Value = '  -0.88%  '; Protect = $false }
    @{ Cell = "D20"; Value = '6.709'; Protect = $true }
    @{ Cell = "E20"; Value = '  -1.15%  '; Protect = $false }
    @{ Cell = "D21"; Value = '17.33'; Protect = $true }
    @{ Cell = "E21"; Value = '  -2.00%  '; Protect = $false }
    @{ Cell = "E22"; Value = '  -0.10%  '; Protect = $false }
    @{ Cell = "D23"; Value = '12.32'; Protect = $true }
    @{ Cell = "E23"; Value = '  -4.02%  '; Protect = $false }
    @{ Cell = "D24"; Value = '23.416.02'; Protect = $false }
    @{ Cell = "E24"; Value = '  -1.18%  '; Protect = $false }
    @{ Cell = "D25"; Value = '2.469'; Protect = $true }
    @{ Cell = "E25"; Value = '  -2.23%  '; Protect = $false }
    @{ Cell = "D26"; Value = '2.891'; Protect = $true }
    @{ Cell = "E26"; Value = '  -4.97%  '; Protect = $false }
    @{ Cell = "D27"; Value = '20.83'; Protect = $true }
    @{ Cell = "E27"; Value = '  -2.08%  '; Protect = $false }
    @{ Cell = "D28"; Value = '152.14'; Protect = $true }
    @{ Cell = "E28"; Value = '  +0.49%  '; Protect = $false }
    @{ Cell = "D29"; Value = '5.194'; Protect = $true }
    @{ Cell = "E29"; Value = '  -0.57%  '; Protect = $false }
    @{ Cell = "D30"; Value = '132.88'; Protect = $true }
    @{ Cell = "E30"; Value = '  -1.14%  '; Protect = $false }
    @{ Cell = "D31"; Value = '1.819.95'; Protect = $false }
    @{ Cell = "E31"; Value = '  -0.91%  '; Protect = $false }
    @{ Cell = "D32"; Value = '6.837'; Protect = $true }
    @{ Cell = "E32"; Value = '  -0.92%  '; Protect = $false }
    @{ Cell = "D33"; Value = '2.130'; Protect = $true }
    @{ Cell = "E33"; Value = '  -4.38%  '; Protect = $false }
    @{ Cell = "D34"; Value = '11.34'; Protect = $true }
    @{ Cell = "E34"; Value = '  -2.82%  '; Protect = $false }
    @{ Cell = "D35"; Value = '0.9755'; Protect = $true }
    @{ Cell = "E35"; Value = '  -8.61%  '; Protect = $false }
    @{ Cell = "D36"; Value = '0.02706'; Protect = $true }
    @{ Cell = "E36"; Value = '  -3.48%  '; Protect = $false }
    @{ Cell = "D37"; Value = '0.08728'; Protect = $true }
    @{ Cell = "E37"; Value = '  -1.12%  '; Protect = $false }
    @{ Cell = "D38"; Value = '0.2423'; Protect = $true }
    @{ Cell = "E38"; Value = '  -3.48%  '; Protect = $false }
    @{ Cell = "D39"; Value = '5.884'; Protect = $true }
    @{ Cell = "E39"; Value = '  -3.34%  '; Protect = $false }
    @{ Cell = "D40"; Value = '13.01'; Protect = $true }
    @{ Cell = "E40"; Value = '  +0.19%  '; Protect = $false }
    @{ Cell = "D41"; Value = '0.06802'; Protect = $true }
    @{ Cell = "E41"; Value = '  -3.87%  '; Protect = $false }
    @{ Cell = "D42"; Value = '0.6836'; Protect = $true }
    @{ Cell = "E42"; Value = '  -2.58%  '; Protect = $false }
    @{ Cell = "D43"; Value = '1.300'; Protect = $true }
    @{ Cell = "E43"; Value = '  -2.91%  '; Protect = $false }
    @{ Cell = "D44"; Value = '15.60'; Protect = $true }
    @{ Cell = "E44"; Value = '  -3.07%  '; Protect = $false }
    @{ Cell = "D45"; Value = '1.000'; Protect = $true }
    @{ Cell = "E45"; Value = '  +0.00%  '; Protect = $false }
    @{ Cell = "D46"; Value = '0.6315'; Protect = $true }
    @{ Cell = "E46"; Value = '  -3.00%  '; Protect = $false }
    @{ Cell = "D47"; Value = '2.240'; Protect = $true }
    @{ Cell = "E47"; Value = '  -3.35%  '; Protect = $false }
    @{ Cell = "D48"; Value = '3.902'; Protect = $true }
    @{ Cell = "E48"; Value = '  -1.30%  '; Protect = $false }
    @{ Cell = "D49"; Value = '0.07702'; Protect = $true }
    @{ Cell = "E49"; Value = '  -3.24%  '; Protect = $false }
    @{ Cell = "D50"; Value = '126.77'; Protect = $true }
    @{ Cell = "E50"; Value = '  -0.87%  '; Protect = $false }
    @{ Cell = "B51"; Value = 'EOS'; Protect = $false }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'; Protect = $false }
    @{ Cell = "D51"; Value = '1.221'; Protect = $true }
    @{ Cell = "E51"; Value = '  +2.04%  '; Protect = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Protect) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}

Write-Output "Applied $($updates.Count) cell updates"
